$d = $word.ActiveDocument

# --- Change 1: "Wartość brutto " -> "Cena" + bookmark(_GoBack) + " brutto " ---
$rng = $d.Content
$found = $rng.Find.Execute("Wartość", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Wartość' text"
}
$wStart = $rng.Start
$wEnd = $rng.End
$bmRange = $d.Range($wEnd, $wEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$leftRng = $d.Range($wStart, $wEnd)
$leftRng.Text = "Cena"

# --- Change 2: merge "{" + "sellerIban" + "}" runs into a single "{sellerIban}" run ---
$d.Content.Find.Execute("{sellerIban}", $true, $false, $false, $false, $false, $true, 1, $false, "{sellerIban}", 2)

Write-Host "Edits applied"
